# Swap the two "requisito" entries so that the LOM3202 (Circuitos Elétricos)
# requirement entry precedes the LOM3221 (Laboratório de Eletrônica) entry
# in the shared strings table. Since sheet1.xml keeps referencing the same
# cells (B24/C24 and B25/C25), achieving this reordering is equivalent to
# swapping the text values shown in those two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lom3202 = "LOM3202 -  Circuitos Elétricos  (Requisito)`n"
$lom3221 = "LOM3221 -  Laboratório de Eletrônica  (Indicação de Conjunto)`n"

$ws.Range("B24").Value = $lom3202
$ws.Range("C24").Value = $lom3202

$ws.Range("B25").Value = $lom3221
$ws.Range("C25").Value = $lom3221
